$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# O11: 572880.59 -> 624605.1
$ws.Range("O11").Value = 624605.1

# O14: 1832.18 -> 2073.75
$ws.Range("O14").Value = 2073.75

# N16: (empty) -> 701.99
$ws.Range("N16").Value = 701.99

# O16: (empty) -> 701.99
$ws.Range("O16").Value = 701.99

# N28: 57376.24 -> 63376.24
$ws.Range("N28").Value = 63376.24
